# Applies the row-level refresh of cryptos.xlsx: updated Price (D) / Volume(1h) (E)
# figures, plus three coin rows (32/33 and 35/36) that swapped rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.493.02'
$ws.Range("E2").Value = '  -3.10%  '

$ws.Range("D3").Value = '2.723.33'
$ws.Range("E3").Value = '  -5.72%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''502.87'
$ws.Range("E5").Value = '  -4.38%  '

$ws.Range("D6").Value = '''140.44'
$ws.Range("E6").Value = '  -1.07%  '

$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("E8").Value = '  -3.86%  '

$ws.Range("D9").Value = '2.737.05'
$ws.Range("E9").Value = '  -5.32%  '

$ws.Range("D10").Value = '''6.05'
$ws.Range("E10").Value = '  +2.33%  '

$ws.Range("D11").Value = '''0.104'
$ws.Range("E11").Value = '  -2.31%  '

$ws.Range("D12").Value = '''0.346'
$ws.Range("E12").Value = '  -2.98%  '

$ws.Range("E13").Value = '  +1.13%  '

$ws.Range("D14").Value = '3.204.95'
$ws.Range("E14").Value = '  -5.44%  '

$ws.Range("D15").Value = '58.638.13'
$ws.Range("E15").Value = '  -2.89%  '

$ws.Range("D16").Value = '''21.62'
$ws.Range("E16").Value = '  -3.99%  '

$ws.Range("D17").Value = '2.729.72'
$ws.Range("E17").Value = '  -5.60%  '

$ws.Range("E18").Value = '  -4.05%  '

$ws.Range("D19").Value = '''4.75'
$ws.Range("E19").Value = '  -3.59%  '

$ws.Range("D20").Value = '''10.95'
$ws.Range("E20").Value = '  -5.10%  '

$ws.Range("D21").Value = '''342.74'
$ws.Range("E21").Value = '  -5.36%  '

$ws.Range("E22").Value = '  -4.20%  '

$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("E24").Value = '  -0.40%  '

$ws.Range("D25").Value = '''62.90'
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("E26").Value = '  -5.07%  '

$ws.Range("D27").Value = '''0.173'
$ws.Range("E27").Value = '  -5.15%  '

$ws.Range("D28").Value = '''0.996'
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").Value = '''7.50'
$ws.Range("E29").Value = '  -3.66%  '

$ws.Range("D30").Value = '0.0₃0826'
$ws.Range("E30").Value = '  -3.33%  '

$ws.Range("D31").Value = '''0.999'
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.60'
$ws.Range("E32").Value = '  -4.19%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''19.12'
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("D34").Value = '''152.13'
$ws.Range("E34").Value = '  +2.78%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '''4.19'
$ws.Range("E35").Value = '  -3.09%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '''5.42'
$ws.Range("E36").Value = '  -2.38%  '

$ws.Range("D37").Value = '''0.947'
$ws.Range("E37").Value = '  -4.48%  '

$ws.Range("E38").Value = '  -5.99%  '

$ws.Range("D39").Value = '''35.79'
$ws.Range("E39").Value = '  -5.17%  '

$ws.Range("E40").Value = '  -7.02%  '

$ws.Range("D41").Value = '''3.54'
$ws.Range("E41").Value = '  -3.22%  '

$ws.Range("D42").Value = '2.189.16'
$ws.Range("E42").Value = '  -5.88%  '

$ws.Range("D43").Value = '''0.0558'
$ws.Range("E43").Value = '  -2.25%  '

$ws.Range("D44").Value = '''0.997'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").Value = '''0.603'
$ws.Range("E45").Value = '  -5.75%  '

$ws.Range("D46").Value = '''18.97'
$ws.Range("E46").Value = '  -7.98%  '

$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("D48").Value = '''4.74'
$ws.Range("E48").Value = '  -6.26%  '

$ws.Range("E49").Value = '  -3.14%  '

$ws.Range("D50").Value = '''0.0885'
$ws.Range("E50").Value = '  -4.67%  '
